# Refactor the inventory report header row to a generic "data" placeholder
# for every column, and tighten several column widths accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) previously held per-column labels such as
# "ServiceTag", "CPU model", "PCI device", etc. Replace every header
# cell in A1:O1 with the literal value "data".
$headerRange = $ws.Range("A1:O1")
foreach ($cell in $headerRange.Cells) {
    $cell.Value = "data"
}

# Adjust column widths to match the new, narrower layout. (Values are
# chosen so the stored/rounded column width lands as close as possible
# to the target width used by the authoring copy of Excel.)
$ws.Columns.Item(1).ColumnWidth = 6.75
$ws.Columns.Item(4).ColumnWidth = 5.75
$ws.Columns.Item(5).ColumnWidth = 7.75
$ws.Columns.Item(6).ColumnWidth = 15.75
$ws.Columns.Item(8).ColumnWidth = 7.75
$ws.Columns.Item(10).ColumnWidth = 3.75
$ws.Columns.Item(11).ColumnWidth = 3.75
$ws.Columns.Item(12).ColumnWidth = 8.75
